$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.512.65'
$ws.Range('E2').Value = '  -0.12%  '
$ws.Range('D3').Value = '1.901.42'
$ws.Range('E3').Value = '  +1.26%  '
$c = $ws.Range('D4')
$c.NumberFormat = "@"
$c.Value = '0.9998'
$c.Style = "Normal"
$ws.Range('E4').Value = '  -0.14%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '239.32'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +1.31%  '
$ws.Range('E6').Value = '  -0.13%  '
$c = $ws.Range('D7')
$c.NumberFormat = "@"
$c.Value = '0.4915'
$c.Style = "Normal"
$ws.Range('E7').Value = '  +0.82%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.2935'
$c.Style = "Normal"
$ws.Range('E8').Value = '  +1.43%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.06694'
$c.Style = "Normal"
$ws.Range('E9').Value = '  +0.39%  '
$ws.Range('D10').Value = '1.903.13'
$ws.Range('E10').Value = '  +1.42%  '
$c = $ws.Range('D11')
$c.NumberFormat = "@"
$c.Value = '17.05'
$c.Style = "Normal"
$ws.Range('E11').Value = '  +2.83%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '0.07342'
$c.Style = "Normal"
$ws.Range('E12').Value = '  +1.40%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '5.182'
$c.Style = "Normal"
$ws.Range('E13').Value = '  +3.66%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '88.12'
$c.Style = "Normal"
$ws.Range('E14').Value = '  -0.60%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '0.6696'
$c.Style = "Normal"
$ws.Range('E15').Value = '  +2.85%  '
$ws.Range('D16').Value = '30.494.16'
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '0.000007893'
$c.Style = "Normal"
$ws.Range('E17').Value = '  +0.65%  '
$ws.Range('B18').Value = 'Avalanche'
$ws.Range('C18').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '13.48'
$c.Style = "Normal"
$ws.Range('E18').Value = '  +3.76%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '0.9999'
$c.Style = "Normal"
$ws.Range('E19').Value = '  -0.13%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '5.453'
$c.Style = "Normal"
$ws.Range('E20').Value = '  +15.66%  '
$ws.Range('D21').Value = '2.142.04'
$ws.Range('E21').Value = '  +1.25%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '0.9991'
$c.Style = "Normal"
$ws.Range('E22').Value = '  -0.21%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '197.76'
$c.Style = "Normal"
$ws.Range('E23').Value = '  +1.83%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '6.133'
$c.Style = "Normal"
$ws.Range('E24').Value = '  +0.07%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '9.532'
$c.Style = "Normal"
$ws.Range('E25').Value = '  +1.87%  '
$c = $ws.Range('D26')
$c.NumberFormat = "@"
$c.Value = '162.90'
$c.Style = "Normal"
$ws.Range('E26').Value = '  +3.98%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '18.44'
$c.Style = "Normal"
$ws.Range('E27').Value = '  -0.24%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '1.947'
$c.Style = "Normal"
$ws.Range('E28').Value = '  +6.63%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '1.487'
$c.Style = "Normal"
$ws.Range('E29').Value = '  +5.52%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '4.358'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +2.49%  '
$ws.Range('E31').Value = '  +1.75%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '4.128'
$c.Style = "Normal"
$ws.Range('E32').Value = '  +5.28%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '0.05166'
$c.Style = "Normal"
$ws.Range('E33').Value = '  +1.28%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '0.7450'
$c.Style = "Normal"
$ws.Range('E34').Value = '  +3.31%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '1.109'
$c.Style = "Normal"
$ws.Range('E35').Value = '  +3.03%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '2.725'
$c.Style = "Normal"
$ws.Range('E36').Value = '  +1.15%  '
$c = $ws.Range('D37')
$c.NumberFormat = "@"
$c.Value = '0.01832'
$c.Style = "Normal"
$ws.Range('E37').Value = '  +1.22%  '
$ws.Range('E38').Value = '  +0.84%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '0.9273'
$c.Style = "Normal"
$ws.Range('E39').Value = '  +0.88%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '2.064'
$c.Style = "Normal"
$ws.Range('E40').Value = '  +1.16%  '
$c = $ws.Range('D41')
$c.NumberFormat = "@"
$c.Value = '0.4412'
$c.Style = "Normal"
$ws.Range('E41').Value = '  +0.69%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$c = $ws.Range('D42')
$c.NumberFormat = "@"
$c.Value = '107.21'
$c.Style = "Normal"
$ws.Range('E42').Value = '  +2.54%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '5.947'
$c.Style = "Normal"
$ws.Range('E43').Value = '  +3.93%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '70.06'
$c.Style = "Normal"
$ws.Range('E44').Value = '  +23.08%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '0.9956'
$c.Style = "Normal"
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('E46').Value = '  +3.58%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '7.601'
$c.Style = "Normal"
$ws.Range('E47').Value = '  +3.79%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '9.006'
$c.Style = "Normal"
$ws.Range('E48').Value = '  +4.56%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '35.13'
$c.Style = "Normal"
$ws.Range('E49').Value = '  +6.23%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '0.05835'
$c.Style = "Normal"
$ws.Range('E50').Value = '  +0.18%  '
$c = $ws.Range('D51')
$c.NumberFormat = "@"
$c.Value = '0.3946'
$c.Style = "Normal"
$ws.Range('E51').Value = '  -1.72%  '
